$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.827.69'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '1.643.89'
$ws.Range("E4").Value = '  -0.44%  '
$ws.Range("D5").Value = '''217.12'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '''0.504'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = '''0.253'
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("D9").Value = '''0.0620'
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").Value = '''19.68'
$ws.Range("E10").Value = '  +3.11%  '
$ws.Range("D12").Value = '1.873.27'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '1.637.93'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '''0.527'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '''66.19'
$ws.Range("E16").Value = '  +2.63%  '
$ws.Range("D17").Value = '26.858.73'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").Value = '''218.06'
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("E21").Value = '  +7.36%  '
$ws.Range("D22").Value = '''4.39'
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("D23").Value = '''2.47'
$ws.Range("E23").Value = '  +6.57%  '
$ws.Range("D24").Value = '''9.16'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("D25").Value = '''146.12'
$ws.Range("E26").Value = '  -0.46%  '
$ws.Range("D27").Value = '''7.33'
$ws.Range("E27").Value = '  +3.33%  '
$ws.Range("E28").Value = '  +0.42%  '
$ws.Range("D29").Value = '''15.83'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").Value = '''0.0511'
$ws.Range("E30").Value = '  +1.69%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '''3.38'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '''2.99'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").Value = '''1.55'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").Value = '1.245.41'
$ws.Range("E36").Value = '  -1.91%  '
$ws.Range("D37").Value = '''0.0175'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("D39").Value = '''0.832'
$ws.Range("E39").Value = '  +3.45%  '
$ws.Range("E40").Value = '  -0.44%  '
$ws.Range("D41").Value = '''0.807'
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").Value = '''5.35'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").Value = '1.785.13'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  -4.26%  '
$ws.Range("D45").Value = '''60.90'
$ws.Range("E45").Value = '  +1.27%  '
$ws.Range("D46").Value = '''91.61'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = '''0.0971'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").Value = '''7.56'
$ws.Range("E51").Value = '  +0.44%  '
